# Rename the inline logo pictures living in the document's headers and
# footers:
#   - Pearson Edexcel logo (footers, first-page + default): image1.png -> image2.png
#   - BTec logo            (headers, first-page + default): image2.jpg -> image1.jpg
#
# InlineShape.Name maps onto the drawing's <wp:docPr name="..."/>. Shapes are
# matched by their (unchanged) AlternativeText/description rather than their
# current Name, since Name only round-trips once the shape has actually been
# selected first (an engine quirk for header/footer-anchored pictures).

$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {

    # Footers: Pearson Edexcel logo -> image2.png
    foreach ($ftr in $sec.Footers) {
        if ($ftr.Exists) {
            foreach ($shp in $ftr.Range.InlineShapes) {
                if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp.Range.Select()
                    $word.Selection.InlineShapes.Item(1).Name = "image2.png"
                }
            }
        }
    }

    # Headers: BTec logo -> image1.jpg
    foreach ($hdr in $sec.Headers) {
        if ($hdr.Exists) {
            foreach ($shp in $hdr.Range.InlineShapes) {
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Range.Select()
                    $word.Selection.InlineShapes.Item(1).Name = "image1.jpg"
                }
            }
        }
    }
}
